$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a far-away staging cell to build text values via formula evaluation
# (so the date-like strings are not auto-converted to date serials/styles),
# then paste-special "values only" into the real target cells so no new
# number-format style gets attached to them. Finally clear the staging cell.
$staging = $ws.Cells.Item(1000, 26)

function Set-TextValue($row, $col, $text) {
    $staging.Formula = '="' + $text + '"'
    $staging.Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4163)  # xlPasteValues
}

# Row 173: 07-09-2021
Set-TextValue 173 1 "07-09-2021"
$ws.Cells.Item(173, 2).Value = 1.48
$ws.Cells.Item(173, 3).Value = 1.91
$ws.Cells.Item(173, 4).Value = 2.37
$ws.Cells.Item(173, 5).Value = 3.04
$ws.Cells.Item(173, 6).Value = -0.45

# Row 174: 08-09-2021
Set-TextValue 174 1 "08-09-2021"
$ws.Cells.Item(174, 2).Value = 1.48
$ws.Cells.Item(174, 3).Value = 1.96
$ws.Cells.Item(174, 4).Value = 2.45
$ws.Cells.Item(174, 5).Value = 3.09
$ws.Cells.Item(174, 6).Value = -0.49

$staging.Clear()
